# Insert 10 new rows of weekly Kiwi price data at row 491,
# pushing the existing rows 491-558 down to become rows 501-568.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 10 rows before row 491 (shifts existing rows 491-558 down to 501-568)
$insertRange = $ws.Range("A491:T500")
$insertRange.EntireRow.Insert()

# Common (unchanged) values for all these rows
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100101007
$categoria = "Kiwi"
$variedad = "Hayward"

function Set-KiwiRow {
    param($row, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $origen, $precioKg, $kgUnidad)

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-KiwiRow 491 44491 "Especial" 150 10000 10000 10000 "$/bandeja 10 kilos" "Región Metropolitana" 1000 10
Set-KiwiRow 492 44491 "Especial" 150 11000 11000 11000 "$/bandeja 10 kilos" "Región de O'Higgins" 1100 10
Set-KiwiRow 493 44491 "Extra (doble especial)" 50 12000 12000 12000 "$/bandeja 10 kilos" "Región Metropolitana" 1200 10
Set-KiwiRow 494 44491 "Extra (doble especial)" 135 13000 13000 13000 "$/bandeja 10 kilos" "Región de O'Higgins" 1300 10
Set-KiwiRow 495 44491 "Extra (doble especial)" 10 430000 430000 430000 "$/bins (450 kilos)" "Región de O'Higgins" 956 450
Set-KiwiRow 496 44491 "Primera" 280 8000 8000 8000 "$/bandeja 10 kilos" "Región Metropolitana" 800 10
Set-KiwiRow 497 44491 "Primera" 220 9000 9000 9000 "$/bandeja 10 kilos" "Región de O'Higgins" 900 10
Set-KiwiRow 498 44491 "Primera" 15 270000 270000 270000 "$/bins (450 kilos)" "Región de O'Higgins" 600 450
Set-KiwiRow 499 44491 "Segunda" 200 7000 7000 7000 "$/bandeja 10 kilos" "Región de O'Higgins" 700 10
Set-KiwiRow 500 44491 "Segunda" 20 220000 220000 220000 "$/bins (450 kilos)" "Región de O'Higgins" 489 450

# Apply the date/time number format to column D for the newly inserted rows (matches the other D cells)
$ws.Range("D491:D500").NumberFormat = "YYYY-MM-DD HH:MM:SS"
